$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (B1:G1) - shared strings reordered so that
# FFR_A, FFR_LF now come before LF_A, LF_FFR
$ws.Range("B1").Value = "FFR_A"
$ws.Range("C1").Value = "FFR_LF"
$ws.Range("D1").Value = "LF_A"
$ws.Range("E1").Value = "LF_FFR"
$ws.Range("F1").Value = "A_FFR"
$ws.Range("G1").Value = "A_LF"

# Update params row (row 2)
$ws.Range("B2").Value = 16.67744507534299
$ws.Range("C2").Value = 4.049661285073961
$ws.Range("D2").Value = -9.705206550582979
$ws.Range("E2").Value = 0.1909481729148204
$ws.Range("F2").Value = 0.003482784764221892
$ws.Range("G2").Value = -0.04298382828175914

# Update pvalue row (row 3)
$ws.Range("B3").Value = 0.1314940620420444
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.000001311694272088104
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.1314940620420459
$ws.Range("G3").Value = 0.000001311694272088104
